$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 94.666664
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 97
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = -316

$ws.Range("H15").Value = 1256.7656
$ws.Range("I15").Value = 1256.7656
$ws.Range("K15").Value = 3770.2968
$ws.Range("M15").Value = -3601.2968

$ws.Range("H21").Value = 18387.2
$ws.Range("I21").Value = 18229.25
$ws.Range("J21").Value = 19019
$ws.Range("K21").Value = 18229.25
$ws.Range("L21").Value = 19019
$ws.Range("M21").Value = -17761.25
$ws.Range("N21").Value = -19955

$ws.Range("H23").Value = 18387.2
$ws.Range("I23").Value = 18229.25
$ws.Range("J23").Value = 19019
$ws.Range("K23").Value = 18229.25
$ws.Range("L23").Value = 19019
$ws.Range("M23").Value = -17995.25
$ws.Range("N23").Value = -19487

$ws.Range("H28").Value = 404.29413
$ws.Range("I28").Value = 387.25
$ws.Range("J28").Value = 445.2
$ws.Range("K28").Value = 387.25
$ws.Range("L28").Value = 445.2
$ws.Range("M28").Value = 97.75
$ws.Range("N28").Value = -1415.2

$ws.Range("H29").Value = 85
$ws.Range("I29").Value = 85
$ws.Range("K29").Value = 255
$ws.Range("M29").Value = 26

$ws.Range("H38").Value = 1476.0454
$ws.Range("I38").Value = 79.454544
$ws.Range("J38").Value = 2872.6365
$ws.Range("K38").Value = 238.363632
$ws.Range("L38").Value = 8617.9095
$ws.Range("M38").Value = 133.636368
$ws.Range("N38").Value = -9361.9095

$ws.Range("H58").Value = 1246.9286
$ws.Range("I58").Value = 145.875
$ws.Range("J58").Value = 2715
$ws.Range("K58").Value = 437.625
$ws.Range("L58").Value = 8145
$ws.Range("M58").Value = -287.625
$ws.Range("N58").Value = -8445

$ws.Range("H87").Value = 34655
$ws.Range("J87").Value = 38238.5
$ws.Range("L87").Value = 38238.5
$ws.Range("N87").Value = -40734.5

$ws.Range("H90").Value = 34655
$ws.Range("J90").Value = 38238.5
$ws.Range("L90").Value = 114715.5
$ws.Range("N90").Value = -127195.5

$ws.Range("H125").Value = 43930.086
$ws.Range("I125").Value = 143306.14
$ws.Range("J125").Value = 453.0625
$ws.Range("K125").Value = 1289755.26
$ws.Range("L125").Value = 4077.5625
$ws.Range("M125").Value = -1287295.26
$ws.Range("N125").Value = -8997.5625

$ws.Range("H127").Value = 928.48
$ws.Range("I127").Value = 310.9
$ws.Range("J127").Value = 997.1
$ws.Range("K127").Value = 932.6999999999999
$ws.Range("L127").Value = 2991.3
$ws.Range("M127").Value = 4027.3
$ws.Range("N127").Value = -12911.3

$ws.Range("H138").Value = 3255.1265
$ws.Range("I138").Value = 2338.5
$ws.Range("J138").Value = 3452.554
$ws.Range("K138").Value = 7015.5
$ws.Range("L138").Value = 10357.662
$ws.Range("M138").Value = -1875.5
$ws.Range("N138").Value = -20637.662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16466.658
$ws.Range("I32").Value = 17620.943
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 17620.943
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -17333.943
$ws.Range("N32").Value = -3574

$ws.Range("H45").Value = 1254
$ws.Range("I45").Value = 1444.5
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1444.5
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -1067.5
$ws.Range("N45").Value = -1754

$ws.Range("H122").Value = 1777.2307
$ws.Range("I122").Value = 1783.3334
$ws.Range("J122").Value = 1772
$ws.Range("K122").Value = 5350.0002
$ws.Range("L122").Value = 5316
$ws.Range("M122").Value = -2900.0002
$ws.Range("N122").Value = -10216

$ws.Range("H132").Value = 2788.8823
$ws.Range("I132").Value = 1875
$ws.Range("J132").Value = 3601.2222
$ws.Range("K132").Value = 5625
$ws.Range("L132").Value = 10803.6666
$ws.Range("M132").Value = -3095
$ws.Range("N132").Value = -15863.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4197
$ws.Range("I105").Value = 3844.3
$ws.Range("J105").Value = 4902.4
$ws.Range("K105").Value = 3844.3
$ws.Range("L105").Value = 4902.4
$ws.Range("M105").Value = -2097.3
$ws.Range("N105").Value = -8396.4

$ws.Range("H107").Value = 21811.16
$ws.Range("I107").Value = 25718.047
$ws.Range("K107").Value = 25718.047
$ws.Range("M107").Value = -23798.047

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 535.1667
$ws.Range("I10").Value = 535.1667
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 535.1667
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -396.1667
$ws.Range("N10").ClearContents()

$ws.Range("H31").Value = 2294.1135
$ws.Range("I31").Value = 1611.1936
$ws.Range("K31").Value = 1611.1936
$ws.Range("M31").Value = -1316.1936

$ws.Range("H34").Value = 2294.1135
$ws.Range("I34").Value = 1611.1936
$ws.Range("K34").Value = 1611.1936
$ws.Range("M34").Value = -1409.1936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1854.6111
$ws.Range("I5").Value = 2198.3333
$ws.Range("K5").Value = 6594.999899999999
$ws.Range("M5").Value = -6482.999899999999

$ws.Range("H34").Value = 1136.25
$ws.Range("I34").Value = 895
$ws.Range("J34").Value = 1216.6666
$ws.Range("K34").Value = 2685
$ws.Range("L34").Value = 3649.9998
$ws.Range("M34").Value = -2601
$ws.Range("N34").Value = -3817.9998

$ws.Range("H39").Value = 2386.6
$ws.Range("J39").Value = 2386.6
$ws.Range("L39").Value = 7159.799999999999
$ws.Range("N39").Value = -7747.799999999999

$ws.Range("H55").Value = 8686.25
$ws.Range("J55").Value = 9855.714
$ws.Range("L55").Value = 29567.142
$ws.Range("N55").Value = -29921.142

$ws.Range("H113").Value = 556185.6
$ws.Range("I113").Value = 2000566.6
$ws.Range("J113").Value = 654.46155
$ws.Range("K113").Value = 6001699.800000001
$ws.Range("L113").Value = 1963.38465
$ws.Range("M113").Value = -5999529.800000001
$ws.Range("N113").Value = -6303.38465

$ws.Range("H122").Value = 744.5952
$ws.Range("I122").Value = 688.4
$ws.Range("K122").Value = 6195.599999999999
$ws.Range("M122").Value = -3745.599999999999

$ws.Range("H123").Value = 915.4
$ws.Range("I123").Value = 915.4
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 2746.2
$ws.Range("L123").Value = 0
$ws.Range("M123").Value = -296.1999999999998
$ws.Range("N123").ClearContents()

$ws.Range("H135").Value = 1854.6111
$ws.Range("I135").Value = 2198.3333
$ws.Range("K135").Value = 19784.9997
$ws.Range("M135").Value = -17249.9997

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5286.0645
$ws.Range("I132").Value = 4804
$ws.Range("J132").Value = 6162.5454
$ws.Range("K132").Value = 14412
$ws.Range("L132").Value = 18487.6362
$ws.Range("M132").Value = -11882
$ws.Range("N132").Value = -23547.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 26266.666
$ws.Range("J54").Value = 26266.666
$ws.Range("L54").Value = 26266.666
$ws.Range("N54").Value = -27306.666

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31240

$ws.Range("H81").Value = 113647.78
$ws.Range("I81").Value = 169175
$ws.Range("J81").Value = 2593.3333
$ws.Range("K81").Value = 338350
$ws.Range("L81").Value = 5186.6666
$ws.Range("M81").Value = -337289
$ws.Range("N81").Value = -7308.6666

$ws.Range("H84").Value = 113647.78
$ws.Range("I84").Value = 169175
$ws.Range("J84").Value = 2593.3333
$ws.Range("K84").Value = 1691750
$ws.Range("L84").Value = 25933.333
$ws.Range("M84").Value = -1686446
$ws.Range("N84").Value = -36541.333

$ws.Range("H132").Value = 1487.4359
$ws.Range("I132").Value = 1168.28
$ws.Range("J132").Value = 2057.3572
$ws.Range("K132").Value = 3504.84
$ws.Range("L132").Value = 6172.071599999999
$ws.Range("M132").Value = -974.8400000000001
$ws.Range("N132").Value = -11232.0716
